$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert a new "Date" column before the current column B ---
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = 11.140625

# --- 2. Header row ---
$ws.Cells.Item(1, 2).Value = "Date"

# --- 3. Row 2 content ---
$ws.Cells.Item(2, 2).Value = [DateTime]"2012-10-22"
$ws.Cells.Item(2, 2).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(2, 3).Value = "il faut que le résultat obtenu soit proportionnel au temps de l'explosion et que cela soit un entier"

# --- 4. Alignment: whole row 2 gets vertical top ---
$xlTop = -4160
$xlCenter = -4108
$ws.Range("A2:G2").VerticalAlignment = $xlTop
$ws.Cells.Item(2, 3).WrapText = $true
$ws.Cells.Item(2, 2).HorizontalAlignment = $xlCenter
$ws.Cells.Item(2, 6).HorizontalAlignment = $xlCenter

# --- 5. Row height for row 2 ---
$ws.Rows.Item(2).RowHeight = 31.5

# --- 6. Add more helper rows at the bottom (F column only), rows 21-29 ---
$xlPasteFormats = -4122
$ws.Range("F20").Copy()
$ws.Range("F21:F29").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- 7. Selection ---
$ws.Range("C13").Select()
